# Applies the "updated instruction excel sheet" edit described by the diff:
#  - A handful of E-column status flags flip from -1 (or 0) to 1
#  - Several new notes are added in column G (new shared strings)
#  - The old note in G21 is removed
#  - The selected cell moves from E43 to E36 (and the saved scroll position
#    no longer pins topLeftCell to A12)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Column E status flips: -1/0 -> 1
# ---------------------------------------------------------------------------
$eRows = 8,11,12,15,16,17,18,19,20,21,30,31,32,33,34,35,41,43,44,46,47,48,49,54
foreach ($r in $eRows) {
    $ws.Range("E$r").Value = 1
}

# ---------------------------------------------------------------------------
# 2. New notes added in column G
# ---------------------------------------------------------------------------
$ws.Range("G2").Value  = "need to implement add overflow"
$ws.Range("G8").Value  = "might have to check imm value"
$ws.Range("G9").Value  = "read more on opcode"
$ws.Range("G10").Value = "read more on opcode"
$ws.Range("G13").Value = "read more on opcode"
$ws.Range("G14").Value = "read more on opcode"

$memNote = "memory not yet implemented"
$ws.Range("G22").Value = $memNote
$ws.Range("G23").Value = $memNote
$ws.Range("G24").Value = $memNote
$ws.Range("G25").Value = $memNote
$ws.Range("G38").Value = $memNote
$ws.Range("G39").Value = $memNote
$ws.Range("G52").Value = $memNote

$ws.Range("G50").Value = "need to implement sub overflow"

# ---------------------------------------------------------------------------
# 3. Drop the stale "still need to make sure exception break..." note
# ---------------------------------------------------------------------------
$ws.Range("G21").ClearContents()

# ---------------------------------------------------------------------------
# 4. Update the saved selection / scroll position
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("E36").Select()
